# Add "Snake and Ladder Problem" as question #33 to the Question List sheet,
# and give a handful of previously "blank-style" C column cells the same
# centered formatting as their row (style index 1 / center alignment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 26-28 get a new, empty, center-aligned C cell (no value).
$ws.Range("C26").HorizontalAlignment = -4108
$ws.Range("C27").HorizontalAlignment = -4108
$ws.Range("C28").HorizontalAlignment = -4108

# Rows 29-33 already have a value in column C; just add the centered style.
$ws.Range("C29").HorizontalAlignment = -4108
$ws.Range("C30").HorizontalAlignment = -4108
$ws.Range("C31").HorizontalAlignment = -4108
$ws.Range("C32").HorizontalAlignment = -4108
$ws.Range("C33").HorizontalAlignment = -4108

# Create new row 34 by copying the formatting of row 33, then filling values.
$ws.Range("A33:F33").Copy()
$ws.Range("A34:F34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "Snake and Ladder Problem"
$ws.Range("C34").Value = "BFS"
$ws.Range("D34").Value = "Graph"
$ws.Range("E34").Value = "medium"
$ws.Range("F34").Value = "GeeksForGeeks"
$ws.Range("C34").HorizontalAlignment = -4108

# Match the final selection state left behind by the edit.
[void]$ws.Range("C26:C34").Select()
